$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1496.9565
$ws.Range("I28").Value = 2176.7273
$ws.Range("J28").Value = 873.8333
$ws.Range("K28").Value = 2176.7273
$ws.Range("L28").Value = 873.8333
$ws.Range("M28").Value = -1691.7273
$ws.Range("N28").Value = -1843.8333
$ws.Range("H116").Value = 3755.2222
$ws.Range("J116").Value = 4000
$ws.Range("L116").Value = 4000
$ws.Range("N116").Value = -10884
$ws.Range("H125").Value = 2161.25
$ws.Range("I125").Value = 1573.875
$ws.Range("K125").Value = 14164.875
$ws.Range("M125").Value = -11704.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5608.5835
$ws.Range("I32").Value = 3663.943
$ws.Range("K32").Value = 3663.943
$ws.Range("M32").Value = -3376.943
$ws.Range("H45").Value = 53441.95
$ws.Range("I45").Value = 84047.914
$ws.Range("J45").Value = 974.5714
$ws.Range("K45").Value = 84047.914
$ws.Range("L45").Value = 974.5714
$ws.Range("M45").Value = -83670.914
$ws.Range("N45").Value = -1728.5714
$ws.Range("H74").Value = 6980660.5
$ws.Range("I74").Value = 10345457
$ws.Range("K74").Value = 10345457
$ws.Range("M74").Value = -10344583
$ws.Range("H77").Value = 6980660.5
$ws.Range("I77").Value = 10345457
$ws.Range("K77").Value = 51727285
$ws.Range("M77").Value = -51722917
$ws.Range("H97").Value = 512.1818
$ws.Range("I97").Value = 393.7143
$ws.Range("K97").Value = 393.7143
$ws.Range("M97").Value = 102.2857
$ws.Range("H102").Value = 1654.2916
$ws.Range("I102").Value = 1481.3158
$ws.Range("J102").Value = 2311.6
$ws.Range("K102").Value = 1481.3158
$ws.Range("L102").Value = 2311.6
$ws.Range("M102").Value = 140.6841999999999
$ws.Range("N102").Value = -5555.6
$ws.Range("H110").Value = 1099.7333
$ws.Range("I110").Value = 863.9091
$ws.Range("J110").Value = 1748.25
$ws.Range("K110").Value = 863.9091
$ws.Range("L110").Value = 1748.25
$ws.Range("M110").Value = 1181.0909
$ws.Range("N110").Value = -5838.25
$ws.Range("H122").Value = 1426.931
$ws.Range("I122").Value = 1415.7727
$ws.Range("K122").Value = 4247.3181
$ws.Range("M122").Value = -1797.3181
$ws.Range("H132").Value = 44551.543
$ws.Range("I132").Value = 3223.6428
$ws.Range("J132").Value = 102410.6
$ws.Range("K132").Value = 9670.928400000001
$ws.Range("L132").Value = 307231.8
$ws.Range("M132").Value = -7140.928400000001
$ws.Range("N132").Value = -312291.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 537.42426
$ws.Range("I94").Value = 483.7
$ws.Range("J94").Value = 620.0769
$ws.Range("K94").Value = 483.7
$ws.Range("L94").Value = 620.0769
$ws.Range("M94").Value = -32.69999999999999
$ws.Range("N94").Value = -1522.0769
$ws.Range("H99").Value = 1714.6666
$ws.Range("I99").Value = 1837.6
$ws.Range("J99").Value = 1100
$ws.Range("K99").Value = 1837.6
$ws.Range("L99").Value = 1100
$ws.Range("M99").Value = -339.5999999999999
$ws.Range("N99").Value = -4096

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1851.25
$ws.Range("I16").Value = 2354.8462
$ws.Range("K16").Value = 2354.8462
$ws.Range("M16").Value = -2067.8462
$ws.Range("H31").Value = 1601.6731
$ws.Range("I31").Value = 1001.7692
$ws.Range("J31").Value = 2201.577
$ws.Range("K31").Value = 1001.7692
$ws.Range("L31").Value = 2201.577
$ws.Range("M31").Value = -706.7692
$ws.Range("N31").Value = -2791.577
$ws.Range("H34").Value = 1601.6731
$ws.Range("I34").Value = 1001.7692
$ws.Range("J34").Value = 2201.577
$ws.Range("K34").Value = 1001.7692
$ws.Range("L34").Value = 2201.577
$ws.Range("M34").Value = -799.7692
$ws.Range("N34").Value = -2605.577
$ws.Range("H113").Value = 1851.25
$ws.Range("I113").Value = 2354.8462
$ws.Range("K113").Value = 2354.8462
$ws.Range("M113").Value = -184.8462
$ws.Range("H132").Value = 1275.4625
$ws.Range("I132").Value = 1184.7858
$ws.Range("J132").Value = 1910.2
$ws.Range("K132").Value = 3554.3574
$ws.Range("L132").Value = 5730.6
$ws.Range("M132").Value = -1024.3574
$ws.Range("N132").Value = -10790.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1580.5
$ws.Range("I114").Value = 322
$ws.Range("J114").Value = 2629.25
$ws.Range("K114").Value = 966
$ws.Range("L114").Value = 7887.75
$ws.Range("M114").Value = 2288
$ws.Range("N114").Value = -14395.75
$ws.Range("H122").Value = 18520360
$ws.Range("J122").Value = 3958.3333
$ws.Range("L122").Value = 35624.9997
$ws.Range("N122").Value = -40524.9997
$ws.Range("H131").Value = 4279.643
$ws.Range("I131").Value = 5338
$ws.Range("J131").Value = 3691.6667
$ws.Range("K131").Value = 16014
$ws.Range("L131").Value = 11075.0001
$ws.Range("M131").Value = -10974
$ws.Range("N131").Value = -21155.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4769.1665
$ws.Range("I70").Value = 4692.222
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 4692.222
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -4422.222
$ws.Range("N70").Value = -5540
$ws.Range("H73").Value = 4769.1665
$ws.Range("I73").Value = 4692.222
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 4692.222
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -3756.222
$ws.Range("N73").Value = -6872
$ws.Range("H102").Value = 1482.6451
$ws.Range("I102").Value = 1315
$ws.Range("J102").Value = 1787.4546
$ws.Range("K102").Value = 1315
$ws.Range("L102").Value = 1787.4546
$ws.Range("M102").Value = 307
$ws.Range("N102").Value = -5031.4546
$ws.Range("H107").Value = 663.7059
$ws.Range("J107").Value = 872.44446
$ws.Range("L107").Value = 872.44446
$ws.Range("N107").Value = -4712.44446
$ws.Range("H113").Value = 1634.52
$ws.Range("I113").Value = 1187.3
$ws.Range("J113").Value = 1932.6666
$ws.Range("K113").Value = 1187.3
$ws.Range("L113").Value = 1932.6666
$ws.Range("M113").Value = 982.7
$ws.Range("N113").Value = -6272.6666
$ws.Range("H132").Value = 1856.841
$ws.Range("I132").Value = 1882.3
$ws.Range("J132").Value = 1802.2858
$ws.Range("K132").Value = 5646.9
$ws.Range("L132").Value = 5406.857400000001
$ws.Range("M132").Value = -3116.9
$ws.Range("N132").Value = -10466.8574

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 38498.7
$ws.Range("I132").Value = 51470.773
$ws.Range("J132").Value = 2825.5
$ws.Range("K132").Value = 154412.319
$ws.Range("L132").Value = 8476.5
$ws.Range("M132").Value = -151882.319
$ws.Range("N132").Value = -13536.5
$ws.Range("H136").Value = 12190
$ws.Range("I136").Value = 13550
$ws.Range("K136").Value = 40650
$ws.Range("M136").Value = -38100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 30627
$ws.Range("J46").Value = 30627
$ws.Range("L46").Value = 30627
$ws.Range("N46").Value = -31089
$ws.Range("H134").Value = 30627
$ws.Range("J134").Value = 30627
$ws.Range("L134").Value = 91881
$ws.Range("N134").Value = -96951
